$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two obsolete DMS config rows (DMSReturnFolderUrl, DMSDocumentUrl).
# Everything below shifts up by two rows; the table/autofilter/dimension
# refs follow automatically.
$ws.Rows("22:23").Delete()

# DMSHomeUrl: repoint at the new prod home site and drop the wrap/indent
# formatting that was special-cased on this row.
$cell = $ws.Range("B20")
$cell.Value = "https://defra.sharepoint.com/sites/EADMSProdhomesite"
$cell.WrapText = $false
$cell.IndentLevel = 0

# DMSSiteUrl: repoint at the new prod site and drop the wrap formatting.
$cell = $ws.Range("B21")
$cell.Value = "https://defra.sharepoint.com/sites/EADMSProd"
$cell.WrapText = $false
$cell.IndentLevel = 0

# DMSEmailTitle: new naming convention for the email/submission title.
$ws.Range("B22").Value = "Waste Return Correspondence {0} {1} - Email and Submission"

# DMSExcelReturnTitle: new naming convention for the excel return title.
$ws.Range("B24").Value = "Waste Return {0} {1}"

# Leave the view scrolled/selected the way the author's session ended up.
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B25").Select()
